$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.483.28"
$ws.Range("E2").Value = "  +4.08%  "
$ws.Range("D3").Value = "1.818.62"
$ws.Range("E3").Value = "  +5.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3839"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3523"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.237"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07824"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.20%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.609"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.25%  "
$ws.Range("D15").Value = "1.815.43"
$ws.Range("E15").Value = "  +5.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.236"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001122"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06734"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.581"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "27.483.80"
$ws.Range("E24").Value = "  +4.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.459"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.687"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.07%  "
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "2.019.54"
$ws.Range("E30").Value = "  +5.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.366"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.062"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08814"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.652"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7014"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2266"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.57%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06502"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02407"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.297"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6601"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.959"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.187"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07328"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.88%  "
